$d = $word.ActiveDocument
Write-Output $d.Paragraphs.Count
Write-Output $d.Content.Text.Length
